# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Rebuilds the "Sources" block at the bottom of the Summary sheet
# (rows 57-79) with the expanded list of citations, and removes the
# single hyperlink that used to sit on the old row 59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlink (old A59) entirely - the new layout has
# no hyperlinked cells, just plain source/url text.
if ($ws.Hyperlinks.Count() -gt 0) {
    $ws.Hyperlinks.Delete()
}

function Set-SourceCell($row, $text, $bold) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $text
    $cell.Font.Bold = $bold
    $cell.Font.Italic = (-not $bold)
    $cell.Font.Underline = $false
}

Set-SourceCell 57 "Source:" $false
Set-SourceCell 58 "" $false
Set-SourceCell 59 "8165.0 - Counts of Australian Businesses (Downloads)" $false
Set-SourceCell 60 "" $false
Set-SourceCell 61 "http://www.abs.gov.au/AUSSTATS/abs@.nsf/Lookup/8165.0Main+Features1Jun%202008%20to%20Jun%202012?OpenDocument" $false
Set-SourceCell 62 "" $false
Set-SourceCell 63 "Others:" $false
Set-SourceCell 64 "" $false
Set-SourceCell 65 "Key Statistics Australian Small Businesses" $false
Set-SourceCell 66 "" $false
Set-SourceCell 67 "http://workspace.unpan.org/sites/internet/Documents/UNPAN92675.pdf" $false
Set-SourceCell 68 "" $false
Set-SourceCell 69 "ABS, Australian Industry:" $false
Set-SourceCell 70 "" $false
Set-SourceCell 71 "http://abs.gov.au/AUSSTATS/abs@.nsf/Lookup/8155.0Main+Features12011-12?OpenDocument" $false
Set-SourceCell 72 "" $false
Set-SourceCell 73 "SME Association of Australia" $false
Set-SourceCell 74 "" $false
Set-SourceCell 75 "http://www.smeaustralia.asn.au/" $false

Set-SourceCell 78 "ABS" $true
Set-SourceCell 79 "Fair Work Act 2009." $false
